# Fix: prevent hidden columns from being labeled as "changed" (ÄNDERUNG)
# Rows 58-105 were false positives (no real difference between the
# FV2304 columns B:K and the FV2310 columns M:V) caused by comparing a
# hidden column. Remove the erroneous "ÄNDERUNG" flag from column L for
# all of those rows, and restyle the rows that start a new segment group
# (previously mis-styled as "changed" rows) to match the normal
# "group header" look already used elsewhere in the sheet (e.g. row 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that are the first row of a new segment ("Segmentname") group and
# which turned out to have no genuine change -> restyle like row 9
# (A/C:K/M:V -> style class "2", B/M -> style class "3").
$groupHeaderRows = @(58, 62, 69, 74, 77, 82, 86, 89, 94, 98, 103)

foreach ($r in $groupHeaderRows) {
    $ws.Range("A9:K9").Copy()
    $ws.Range("A$r`:K$r").PasteSpecial(-4122)
    $ws.Range("M9:V9").Copy()
    $ws.Range("M$r`:V$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# All affected rows (58-105) lose the false "ÄNDERUNG" label in column L.
$ws.Range("L58:L105").ClearContents()
$ws.Range("L3").Copy()
$ws.Range("L58:L105").PasteSpecial(-4122)
$excel.CutCopyMode = 0
